# Refresh the crypto price/volume table (GitHub Actions data pull).
# For cells whose new text parses as a plain number (e.g. "263.58"), the
# cell's NumberFormat is forced to Text ("@") first so Excel keeps storing
# the value as a string (matching the source data feed) instead of
# silently converting it to a numeric cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.602.16'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '1.851.82'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '263.58'
$ws.Range("E5").Value = '  +2.11%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5266'
$ws.Range("E7").Value = '  +1.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3244'
$ws.Range("E8").Value = '  +0.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06807'
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.97'
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7844'
$ws.Range("E11").Value = '  +2.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07766'
$ws.Range("E12").Value = '  +1.07%  '
$ws.Range("D13").Value = '1.857.02'
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.82'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.042'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '26.630.56'
$ws.Range("E20").Value = '  +1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.642'
$ws.Range("E21").Value = '  +2.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.495'
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.020'
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.88'
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.174'
$ws.Range("E25").Value = '  -5.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.680'
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.02'
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '111.68'
$ws.Range("E28").Value = '  +0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.190'
$ws.Range("E29").Value = '  +0.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.115'
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08723'
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04871'
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7245'
$ws.Range("E33").Value = '  +7.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.874'
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.118'
$ws.Range("E36").Value = '  +0.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.277'
$ws.Range("E37").Value = '  +4.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01792'
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.4872'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9020'
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '111.39'
$ws.Range("E41").Value = '  +0.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.969'
$ws.Range("E42").Value = '  -2.26%  '
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.701'
$ws.Range("E44").Value = '  +0.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4203'
$ws.Range("E45").Value = '  +0.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05880'
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.029'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1241'
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.15'
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.8911'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.09'
$ws.Range("E51").Value = '  +2.12%  '
